$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1,1).Range.Text = "41-23="
$t.Cell(1,2).Range.Text = "84-27="
$t.Cell(1,3).Range.Text = "39+42="
$t.Cell(1,4).Range.Text = "87-18="
$t.Cell(1,5).Range.Text = "97-48="
$t.Cell(2,1).Range.Text = "64-35="
$t.Cell(2,2).Range.Text = "80-33="
$t.Cell(2,3).Range.Text = "15+69="
$t.Cell(2,4).Range.Text = "47+26="
$t.Cell(2,5).Range.Text = "3+69="
$t.Cell(3,1).Range.Text = "17+58="
$t.Cell(3,2).Range.Text = "51-44="
$t.Cell(3,3).Range.Text = "11-6="
$t.Cell(3,4).Range.Text = "71-38="
$t.Cell(3,5).Range.Text = "82-49="
$t.Cell(4,1).Range.Text = "53-7="
$t.Cell(4,2).Range.Text = "36+45="
$t.Cell(4,3).Range.Text = "63-36="
$t.Cell(4,4).Range.Text = "93-69="
$t.Cell(4,5).Range.Text = "91-26="
$t.Cell(5,1).Range.Text = "43-29="
$t.Cell(5,2).Range.Text = "29+32="
$t.Cell(5,3).Range.Text = "51-3="
$t.Cell(5,4).Range.Text = "60-9="
$t.Cell(5,5).Range.Text = "23-6="
$t.Cell(6,1).Range.Text = "5+36="
$t.Cell(6,2).Range.Text = "16+29="
$t.Cell(6,3).Range.Text = "43-28="
$t.Cell(6,4).Range.Text = "29+63="
$t.Cell(6,5).Range.Text = "81-43="
$t.Cell(7,1).Range.Text = "38+25="
$t.Cell(7,2).Range.Text = "75-16="
$t.Cell(7,3).Range.Text = "17+39="
$t.Cell(7,4).Range.Text = "18+47="
$t.Cell(7,5).Range.Text = "90-29="
$t.Cell(8,1).Range.Text = "37+18="
$t.Cell(8,2).Range.Text = "55+7="
$t.Cell(8,3).Range.Text = "29+4="
$t.Cell(8,4).Range.Text = "84-36="
$t.Cell(8,5).Range.Text = "50-43="
$t.Cell(9,1).Range.Text = "56-27="
$t.Cell(9,2).Range.Text = "61-54="
$t.Cell(9,3).Range.Text = "29+26="
$t.Cell(9,4).Range.Text = "9+22="
$t.Cell(9,5).Range.Text = "7+75="
$t.Cell(10,1).Range.Text = "16+9="
$t.Cell(10,2).Range.Text = "53-4="
$t.Cell(10,3).Range.Text = "82-34="
$t.Cell(10,4).Range.Text = "20-6="
$t.Cell(10,5).Range.Text = "82-16="
$t.Cell(11,1).Range.Text = "26+27="
$t.Cell(11,2).Range.Text = "66+7="
$t.Cell(11,3).Range.Text = "77+18="
$t.Cell(11,4).Range.Text = "28+59="
$t.Cell(11,5).Range.Text = "85-79="
$t.Cell(12,1).Range.Text = "24+9="
$t.Cell(12,2).Range.Text = "42-38="
$t.Cell(12,3).Range.Text = "49+36="
$t.Cell(12,4).Range.Text = "74-38="
$t.Cell(12,5).Range.Text = "60-7="
$t.Cell(13,1).Range.Text = "15+36="
$t.Cell(13,2).Range.Text = "6+45="
$t.Cell(13,3).Range.Text = "66+7="
$t.Cell(13,4).Range.Text = "76+18="
$t.Cell(13,5).Range.Text = "19+73="
$t.Cell(14,1).Range.Text = "46-27="
$t.Cell(14,2).Range.Text = "97-29="
$t.Cell(14,3).Range.Text = "81-32="
$t.Cell(14,4).Range.Text = "80-17="
$t.Cell(14,5).Range.Text = "61-5="
$t.Cell(15,1).Range.Text = "87+8="
$t.Cell(15,2).Range.Text = "15+37="
$t.Cell(15,3).Range.Text = "7+29="
$t.Cell(15,4).Range.Text = "90-73="
$t.Cell(15,5).Range.Text = "39+26="
$t.Cell(16,1).Range.Text = "15+56="
$t.Cell(16,2).Range.Text = "70-2="
$t.Cell(16,3).Range.Text = "37+37="
$t.Cell(16,4).Range.Text = "8+56="
$t.Cell(16,5).Range.Text = "48+4="
$t.Cell(17,1).Range.Text = "40-29="
$t.Cell(17,2).Range.Text = "5+79="
$t.Cell(17,3).Range.Text = "75-29="
$t.Cell(17,4).Range.Text = "90-72="
$t.Cell(17,5).Range.Text = "15+67="
$t.Cell(18,1).Range.Text = "12-9="
$t.Cell(18,2).Range.Text = "46+6="
$t.Cell(18,3).Range.Text = "7+44="
$t.Cell(18,4).Range.Text = "96-18="
$t.Cell(18,5).Range.Text = "77-49="
$t.Cell(19,1).Range.Text = "89+9="
$t.Cell(19,2).Range.Text = "42+49="
$t.Cell(19,3).Range.Text = "66+7="
$t.Cell(19,4).Range.Text = "52-8="
$t.Cell(19,5).Range.Text = "64-29="
$t.Cell(20,1).Range.Text = "56+7="
$t.Cell(20,2).Range.Text = "30-12="
$t.Cell(20,3).Range.Text = "17+44="
$t.Cell(20,4).Range.Text = "83-15="
$t.Cell(20,5).Range.Text = "72-53="
